# fix off by one error
# - "Texture Rotation API" and "Multitexture API" tasks were completed/removed,
#   and "Fix remaining consistency issues" (now directly after "Fix button text
#   positioning") is marked as Completed ("Yes").
# - One extra blank separator row between "Input Method API" and the
#   "Lightning 2.0 Pre-Final Tasks" section header is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stash the formatting of the lone stray empty formatted cell at C19 so we can
# put it back in the same visual spot once the rows above it shift around.
$ws.Range("C19").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats

# Remove the two rows for "Texture Rotation API" (row 9) and "Multitexture API"
# (row 10) entirely - they shift everything below up by two rows.
$ws.Rows("9:10").Delete()

# The row that used to be "Fix remaining consistency issues" (row 11) is now
# row 9; mark it completed same as the rows above it (copy formatting - the
# green "Yes" fill - from the cell directly above it, then set the value).
$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B9").Value = "Yes"
$excel.CutCopyMode = $false

# Remove one of the two now-redundant blank rows separating "Input Method API"
# (now row 23) from the "Lightning 2.0 Pre-Final Tasks" header (now row 26).
$ws.Rows("24:24").Delete()

# The row delete above carried the stray formatted cell from C19 down to C17;
# put it back where it started (C19) and restore C17 to the default look.
$ws.Range("C17").Clear()
$ws.Range("E1").Copy()
$ws.Range("C19").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E1").Clear()
$excel.CutCopyMode = $false

# Update the active selection to match the saved state of the workbook.
$ws.Range("A11").Select()
